$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
# Row 112
$ws.Range("H112").Value = 2014.2632
$ws.Range("J112").Value = 2230.6875
$ws.Range("L112").Value = 6692.0625
$ws.Range("N112").Value = -8908.0625
# Row 121
$ws.Range("H121").Value = 3950.3333
$ws.Range("J121").Value = 3950.3333
$ws.Range("L121").Value = 11850.9999
$ws.Range("N121").Value = -15344.9999
# Row 131
$ws.Range("H131").Value = 65549
$ws.Range("I131").Value = 93050.55
$ws.Range("J131").Value = 5045.6
$ws.Range("K131").Value = 279151.65
$ws.Range("L131").Value = 15136.8
$ws.Range("M131").Value = -274111.65
$ws.Range("N131").Value = -25216.8
# Row 132
$ws.Range("H132").Value = 2930.15
$ws.Range("I132").Value = 2428.743
$ws.Range("J132").Value = 6440
$ws.Range("K132").Value = 7286.228999999999
$ws.Range("L132").Value = 19320
$ws.Range("M132").Value = -4756.228999999999
$ws.Range("N132").Value = -24380
# Row 141
$ws.Range("H141").Value = 501.5366
$ws.Range("I141").Value = 501.5366
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 1504.6098
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 3675.3902
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets("ARM")
# Row 32
$ws.Range("H32").Value = 10313414
$ws.Range("I32").Value = 12660704
$ws.Range("J32").Value = 11421.167
$ws.Range("K32").Value = 12660704
$ws.Range("L32").Value = 11421.167
$ws.Range("M32").Value = -12660417
$ws.Range("N32").Value = -11995.167
# Row 102
$ws.Range("H102").Value = 1286.4546
$ws.Range("I102").Value = 1340
$ws.Range("J102").Value = 1241.8334
$ws.Range("K102").Value = 1340
$ws.Range("L102").Value = 1241.8334
$ws.Range("M102").Value = 282
$ws.Range("N102").Value = -4485.8334
# Row 139
$ws.Range("H139").Value = 27273
$ws.Range("J139").Value = 29428.75
$ws.Range("L139").Value = 29428.75
$ws.Range("N139").Value = -39708.75

$ws = $wb.Worksheets("BSM")
# Row 53
$ws.Range("H53").Value = 30000
$ws.Range("J53").Value = 30000
$ws.Range("L53").Value = 30000
$ws.Range("N53").Value = -31148
# Row 94
$ws.Range("H94").Value = 428
$ws.Range("I94").Value = 377
$ws.Range("K94").Value = 377
$ws.Range("M94").Value = 74
# Row 99
$ws.Range("H99").Value = 1979
$ws.Range("I99").Value = 774.2857
$ws.Range("K99").Value = 774.2857
$ws.Range("M99").Value = 723.7143

$ws = $wb.Worksheets("CRP")
# Row 31
$ws.Range("H31").Value = 3367.9744
$ws.Range("I31").Value = 2171.5881
$ws.Range("J31").Value = 4292.4546
$ws.Range("K31").Value = 2171.5881
$ws.Range("L31").Value = 4292.4546
$ws.Range("M31").Value = -1876.5881
$ws.Range("N31").Value = -4882.4546
# Row 34
$ws.Range("H34").Value = 3367.9744
$ws.Range("I34").Value = 2171.5881
$ws.Range("J34").Value = 4292.4546
$ws.Range("K34").Value = 2171.5881
$ws.Range("L34").Value = 4292.4546
$ws.Range("M34").Value = -1969.5881
$ws.Range("N34").Value = -4696.4546
# Row 122
$ws.Range("H122").Value = 2091.4
$ws.Range("I122").Value = 1400
$ws.Range("J122").Value = 2782.8
$ws.Range("K122").Value = 4200
$ws.Range("L122").Value = 8348.400000000001
$ws.Range("M122").Value = -1750
$ws.Range("N122").Value = -13248.4

$ws = $wb.Worksheets("CUL")
# Row 5
$ws.Range("H5").Value = 539.7353000000001
$ws.Range("I5").Value = 474.2414
$ws.Range("K5").Value = 1422.7242
$ws.Range("M5").Value = -1310.7242
# Row 80
$ws.Range("H80").Value = 2488.889
$ws.Range("J80").Value = 2488.889
$ws.Range("L80").Value = 7466.667
$ws.Range("N80").Value = -9338.667000000001
# Row 83
$ws.Range("H83").Value = 2488.889
$ws.Range("J83").Value = 2488.889
$ws.Range("L83").Value = 22400.001
$ws.Range("N83").Value = -31760.001
# Row 135
$ws.Range("H135").Value = 539.7353000000001
$ws.Range("I135").Value = 474.2414
$ws.Range("K135").Value = 4268.1726
$ws.Range("M135").Value = -1733.1726

$ws = $wb.Worksheets("GSM")
# Row 97
$ws.Range("H97").Value = 1136.3334
$ws.Range("I97").Value = 1136.3334
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1136.3334
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -640.3334
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets("LTW")
# Row 93
$ws.Range("H93").Value = 10383
$ws.Range("I93").Value = 21341
$ws.Range("J93").Value = 1251.3334
$ws.Range("K93").Value = 21341
$ws.Range("L93").Value = 1251.3334
$ws.Range("M93").Value = -20093
$ws.Range("N93").Value = -3747.3334
# Row 100
$ws.Range("H100").Value = 102000320
$ws.Range("I100").Value = 3333535.2
$ws.Range("K100").Value = 3333535.2
$ws.Range("M100").Value = -3332994.2

$ws = $wb.Worksheets("WVR")
# Row 14
$ws.Range("H14").Value = 7666.6665
$ws.Range("I14").Value = 7666.6665
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 7666.6665
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -7498.6665
$ws.Range("N14").ClearContents()
# Row 46
$ws.Range("H46").Value = 38095.332
$ws.Range("J46").Value = 38095.332
$ws.Range("L46").Value = 38095.332
$ws.Range("N46").Value = -38557.332
# Row 96
$ws.Range("H96").Value = 1669
$ws.Range("I96").Value = 1060.5
$ws.Range("J96").Value = 1912.4
$ws.Range("K96").Value = 1060.5
$ws.Range("L96").Value = 1912.4
$ws.Range("M96").Value = 312.5
$ws.Range("N96").Value = -4658.4
# Row 100
$ws.Range("H100").Value = 700.36365
$ws.Range("I100").Value = 572
$ws.Range("J100").Value = 925
$ws.Range("K100").Value = 1144
$ws.Range("L100").Value = 1850
$ws.Range("M100").Value = -603
$ws.Range("N100").Value = -2932
# Row 122
$ws.Range("H122").Value = 1845.4348
$ws.Range("I122").Value = 1649.1428
$ws.Range("J122").Value = 2150.7778
$ws.Range("K122").Value = 4947.428400000001
$ws.Range("L122").Value = 6452.3334
$ws.Range("M122").Value = -2497.428400000001
$ws.Range("N122").Value = -11352.3334
# Row 126
$ws.Range("H126").Value = 2474.9583
$ws.Range("I126").Value = 1165.8125
$ws.Range("J126").Value = 5093.25
$ws.Range("K126").Value = 3497.4375
$ws.Range("L126").Value = 15279.75
$ws.Range("M126").Value = -1027.4375
$ws.Range("N126").Value = -20219.75
# Row 132
$ws.Range("H132").Value = 2209.3655
$ws.Range("I132").Value = 2103.2
$ws.Range("J132").Value = 2427.9412
$ws.Range("K132").Value = 6309.599999999999
$ws.Range("L132").Value = 7283.823600000001
$ws.Range("M132").Value = -3779.599999999999
$ws.Range("N132").Value = -12343.8236
# Row 134
$ws.Range("H134").Value = 38095.332
$ws.Range("J134").Value = 38095.332
$ws.Range("L134").Value = 114285.996
$ws.Range("N134").Value = -119355.996
